$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "cryptos" table with refreshed price / volume data.
# For any new value that Excel would otherwise auto-parse as a number
# (e.g. "1.00", "7.07"), the cell is briefly switched to Text format
# before the assignment and then restored to the default "Normal" style,
# so the written value stays a plain string (matching the source data)
# without leaving a numeric format behind on the cell.

$ws.Range("D2").Value = "58.363.55"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "2.995.32"
$ws.Range("E3").Value = "  +3.65%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  +4.93%  "
$ws.Range("D9").Value = "2.987.86"
$ws.Range("E9").Value = "  +3.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.132"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.31%  "
$ws.Range("E12").Value = "  +5.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.04%  "
$ws.Range("E15").Value = "  +3.23%  "
$ws.Range("D16").Value = "3.495.33"
$ws.Range("E16").Value = "  +3.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.04%  "
$ws.Range("D18").Value = "2.995.70"
$ws.Range("E18").Value = "  +3.70%  "
$ws.Range("D19").Value = "58.312.85"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "423.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.713"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.96%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E28").Value = "  +3.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0978"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.967"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.02%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.47%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0740"
$ws.Range("E36").Value = "  +20.16%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "48.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +17.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "395.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0352"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.96%  "
$ws.Range("E43").Value = "  +2.79%  "
$ws.Range("D44").Value = "2.727.61"
$ws.Range("E44").Value = "  +4.66%  "
$ws.Range("E45").Value = "  +8.20%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.51%  "
$ws.Range("E48").Value = "  +4.69%  "
$ws.Range("E49").Value = "  +2.90%  "
$ws.Range("E50").Value = "  +4.59%  "
$ws.Range("E51").Value = "  +4.84%  "
